$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.744136364917767
$ws.Range("B2").Value = 61.84429565824441
$ws.Range("C2").Value = 57.37533254787634
$ws.Range("D2").Value = 5.730210727470782
$ws.Range("E2").Value = 5.732369348379622
$ws.Range("F2").Value = 3.540140874333925
